$wb = $excel.ActiveWorkbook

# The UK market test sheet is a duplicate of the Poland sheet (same layout),
# inserted right after it.
$poland = $wb.Worksheets.Item("Poland")
$poland.Copy($null, $poland)
$uk = $wb.Worksheets.Item("Poland (2)")
$uk.Name = "UK"

# Fill in the market-specific values (note: set B4 before B2 so the new
# shared strings are appended in the same order as the source workbook).
$uk.Range("B4").Value = "NGC-2741/T3345/T3343/T3342"
$uk.Range("B2").Value = "UK Market"

# Insert an extra "GMPIM" row before the existing "PR1D2" row, copying the
# row above it first so the new row keeps the same cell styling/border.
$uk.Rows.Item(9).Insert()
$uk.Range("A10").Copy($uk.Range("A9"))
$uk.Range("A9").Value = "GMPIM"

$uk.Activate()
$null = $uk.Range("A9").Select()
